# Updated cryptos list with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency rows in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text, matching the workbook's
# original inline-string cells (e.g. "27.137.48", "0.07261", "  -0.28%  ").
# Without this, Excel's COM layer auto-converts plain numeric-looking
# strings (like "0.07261") into real numbers. Temporarily forcing the
# cell's number format to Text ("@") prevents that conversion, and
# resetting the style back to "Normal" afterwards avoids leaving any
# stray formatting behind on the cell.
function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" '27.137.48'
Set-TextValue $ws "E2" '  -0.28%  '
Set-TextValue $ws "E3" '  -0.85%  '
Set-TextValue $ws "E4" '  +0.10%  '
Set-TextValue $ws "D5" '306.83'
Set-TextValue $ws "E6" '  +0.03%  '
Set-TextValue $ws "D7" '0.5221'
Set-TextValue $ws "E7" '  -0.52%  '
Set-TextValue $ws "D8" '0.3757'
Set-TextValue $ws "E8" '  -0.63%  '
Set-TextValue $ws "D9" '0.07261'
Set-TextValue $ws "D10" '21.10'
Set-TextValue $ws "E10" '  -0.69%  '
Set-TextValue $ws "D11" '0.8994'
Set-TextValue $ws "E11" '  +0.05%  '
Set-TextValue $ws "D12" '0.08159'
Set-TextValue $ws "E12" '  +6.05%  '
Set-TextValue $ws "D13" '1.934.40'
Set-TextValue $ws "E13" '  +1.46%  '
Set-TextValue $ws "D14" '96.22'
Set-TextValue $ws "D15" '5.283'
Set-TextValue $ws "E15" '  +0.25%  '
Set-TextValue $ws "E16" '  +0.08%  '
Set-TextValue $ws "D17" '0.000008573'
Set-TextValue $ws "E17" '  -0.94%  '
Set-TextValue $ws "D18" '14.58'
Set-TextValue $ws "E18" '  +0.43%  '
Set-TextValue $ws "E19" '  +0.13%  '
Set-TextValue $ws "D20" '27.157.98'
Set-TextValue $ws "E20" '  -0.45%  '
Set-TextValue $ws "D21" '5.080'
Set-TextValue $ws "E21" '  -0.11%  '
Set-TextValue $ws "E22" '  +0.41%  '
Set-TextValue $ws "D23" '6.399'
Set-TextValue $ws "E23" '  -0.79%  '
Set-TextValue $ws "D24" '147.82'
Set-TextValue $ws "E24" '  +1.46%  '
Set-TextValue $ws "D25" '2.283'
Set-TextValue $ws "E25" '  -1.85%  '
Set-TextValue $ws "D26" '18.16'
Set-TextValue $ws "E26" '  +0.12%  '
Set-TextValue $ws "D27" '1.741'
Set-TextValue $ws "E27" '  +0.05%  '
Set-TextValue $ws "D28" '114.97'
Set-TextValue $ws "E28" '  +0.08%  '
Set-TextValue $ws "D29" '4.783'
Set-TextValue $ws "D30" '4.842'
Set-TextValue $ws "E30" '  -2.65%  '
Set-TextValue $ws "D31" '0.09217'
Set-TextValue $ws "E31" '  -0.21%  '
Set-TextValue $ws "D32" '0.05049'
Set-TextValue $ws "E32" '  -0.62%  '
Set-TextValue $ws "D33" '0.7877'
Set-TextValue $ws "E33" '  -4.12%  '
Set-TextValue $ws "D35" '3.420'
Set-TextValue $ws "E35" '  +3.32%  '
Set-TextValue $ws "D36" '2.969'
Set-TextValue $ws "E36" '  -0.59%  '
Set-TextValue $ws "D37" '2.581'
Set-TextValue $ws "E37" '  -0.93%  '
Set-TextValue $ws "D38" '0.5693'
Set-TextValue $ws "E38" '  +0.32%  '
Set-TextValue $ws "D39" '0.01976'
Set-TextValue $ws "E39" '  -0.75%  '
Set-TextValue $ws "E40" '  -0.28%  '
Set-TextValue $ws "D41" '9.007'
Set-TextValue $ws "E41" '  -0.04%  '
Set-TextValue $ws "D42" '6.554'
Set-TextValue $ws "E42" '  -1.35%  '
Set-TextValue $ws "E43" '  -2.70%  '
Set-TextValue $ws "E44" '  -0.17%  '
Set-TextValue $ws "D45" '0.4855'
Set-TextValue $ws "E45" '  +0.21%  '
Set-TextValue $ws "E46" '  +0.04%  '
Set-TextValue $ws "D47" '10.09'
Set-TextValue $ws "E47" '  -1.64%  '
Set-TextValue $ws "E48" '  -0.30%  '
Set-TextValue $ws "D49" '38.08'
Set-TextValue $ws "E49" '  +1.25%  '
Set-TextValue $ws "D50" '63.42'
Set-TextValue $ws "E50" '  -0.67%  '
Set-TextValue $ws "D51" '0.05934'
Set-TextValue $ws "E51" '  -0.02%  '
